# "Generate Report for Handoff"
#
# The 0b78385b-...md file has finished translation and is now ready for
# handoff, while e0956004-...md is still in translation. The localization
# report is regenerated: on every sheet the two source-file rows are
# re-ordered (e0956004 first, 0b78385b second) and the 0b78385b row's
# status/handoff-datetime columns are refreshed to reflect its new
# "Ready for handoff" state.

$wb = $excel.ActiveWorkbook

$mdUrl0b  = "https://github.com/OpenLocalizationTest/oltest/blob/bd15f9ffc8aa18ab2d943b695fcfce5061e9f7d5/e2e/0b78385b-b68e-4c64-a372-8ebd007500b5.md"
$mdUrlE0  = "https://github.com/OpenLocalizationTest/oltest/blob/bd15f9ffc8aa18ab2d943b695fcfce5061e9f7d5/e2e/e0956004-b38b-46a0-96f5-1652e41ddf97.md"

$mdName0b  = "0b78385b-b68e-4c64-a372-8ebd007500b5.md"
$mdNameE0  = "e0956004-b38b-46a0-96f5-1652e41ddf97.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-13-20 18:13:43"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl0b, "", "", $mdNameE0)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrlE0, "", "", $mdName0b)

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$xlfUrlZhCn0b = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08c1770b78e808eeb389506812c46d1705627f81/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/0b78385b-b68e-4c64-a372-8ebd007500b5.1e606bd3a11a53369d209f9cf42fd811ab43e80f.zh-cn.xlf"
$xlfUrlZhCnE0 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/08c1770b78e808eeb389506812c46d1705627f81/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e0956004-b38b-46a0-96f5-1652e41ddf97.aefdd078ebfa948027c359ff4233c6a58e2eb94d.zh-cn.xlf"

$xlfNameZhCn0b = "0b78385b-b68e-4c64-a372-8ebd007500b5.1e606bd3a11a53369d209f9cf42fd811ab43e80f.zh-cn.xlf"
$xlfNameZhCnE0 = "e0956004-b38b-46a0-96f5-1652e41ddf97.aefdd078ebfa948027c359ff4233c6a58e2eb94d.zh-cn.xlf"

$ws.Range("A2").Value = $mdNameE0
$ws.Range("D2").Value = $xlfNameZhCnE0

$ws.Range("A3").Value = $mdName0b
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = $xlfNameZhCn0b
$ws.Range("E3").Value = "2016-03-20 18:13:40"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl0b, "", "", $mdNameE0)
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl0b, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfUrlZhCn0b, "", "", $xlfNameZhCnE0)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrlE0, "", "", $mdName0b)
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrlE0, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfUrlZhCnE0, "", "", $xlfNameZhCn0b)

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$xlfUrlDeDe0b = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b1434f1683c0709c641997eb015c67affe558834/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/0b78385b-b68e-4c64-a372-8ebd007500b5.1e606bd3a11a53369d209f9cf42fd811ab43e80f.de-de.xlf"
$xlfUrlDeDeE0 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b1434f1683c0709c641997eb015c67affe558834/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e0956004-b38b-46a0-96f5-1652e41ddf97.aefdd078ebfa948027c359ff4233c6a58e2eb94d.de-de.xlf"

$xlfNameDeDe0b = "0b78385b-b68e-4c64-a372-8ebd007500b5.1e606bd3a11a53369d209f9cf42fd811ab43e80f.de-de.xlf"
$xlfNameDeDeE0 = "e0956004-b38b-46a0-96f5-1652e41ddf97.aefdd078ebfa948027c359ff4233c6a58e2eb94d.de-de.xlf"

$ws.Range("A2").Value = $mdNameE0
$ws.Range("D2").Value = $xlfNameDeDeE0

$ws.Range("A3").Value = $mdName0b
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = $xlfNameDeDe0b
$ws.Range("E3").Value = "2016-03-20 18:13:43"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl0b, "", "", $mdNameE0)
$ws.Hyperlinks.Add($ws.Range("B2"), $mdUrl0b, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D2"), $xlfUrlDeDe0b, "", "", $xlfNameDeDeE0)
$ws.Hyperlinks.Add($ws.Range("A3"), $mdUrlE0, "", "", $mdName0b)
$ws.Hyperlinks.Add($ws.Range("B3"), $mdUrlE0, "", "", ".md")
$ws.Hyperlinks.Add($ws.Range("D3"), $xlfUrlDeDeE0, "", "", $xlfNameDeDe0b)
